# Apply crypto price/volume updates to match Sun Jan  7 18:14:20 UTC 2024 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.566.69'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').Value = '2.241.38'
$ws.Range('E3').Value = '  -0.24%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.40'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.65'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.88%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.570'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.66%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.518'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.23'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0802'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.20'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').Value = '2.582.33'
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('D15').Value = '2.237.75'
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.831'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('E17').Value = '  -0.65%  '
$ws.Range('D18').Value = '44.327.56'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('E19').Value = '  -3.52%  '
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.81'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.26%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.18'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.23'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.47'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.31%  '
$ws.Range('E24').Value = '  -6.42%  '
$ws.Range('E25').Value = '  -2.44%  '
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.35'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +6.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.74'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.07'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.89'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.87'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '149.59'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0784'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.09%  '
$ws.Range('E35').Value = '  -3.02%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.108'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.88'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +6.27%  '
$ws.Range('E38').Value = '  -1.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.06'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.99%  '
$ws.Range('E40').Value = '  -7.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.80'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0297'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').Value = '1.810.10'
$ws.Range('E44').Value = '  +2.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.77'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +12.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '81.09'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.188'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '98.29'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.15%  '
$ws.Range('E49').Value = '  -2.84%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.37'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '53.95'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.58%  '
